$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fa03404f3d0>),
                (''model'',
                 XGBClassifier(base_score=None, booster=None, callbacks=None,
                               colsample_bylevel=None, colsample_bynode=None,
                               colsample_bytree=0.8, early_stopping_rounds=None,
                               enable_categorical=False, eval_metric=None,
                               feature_types=None, gamma=0.1, gpu_id=None,
                               grow_policy=None, importance_type=None,
                               interaction_constraints=None, learning_rate=0.01,
                               max_bin=None, max_cat_threshold=None,
                               max_cat_to_onehot=None, max_delta_step=None,
                               max_depth=3, max_leaves=None,
                               min_child_weight=None, missing=nan,
                               monotone_constraints=None, n_estimators=50,
                               n_jobs=None, num_parallel_tree=None,
                               predictor=None, random_state=42, ...))])'
$ws.Range("B2").Value = 0.6781118881118882
$ws.Range("C2").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fa004493f10>, ''scaler'': MinMaxScaler(), ''model__subsample'': 0.5, ''model__n_estimators'': 50, ''model__max_depth'': 3, ''model__learning_rate'': 0.01, ''model__gamma'': 0.1, ''model__colsample_bytree'': 0.8}'
$ws.Range("D2").Value = 0.9780858117341675
$ws.Range("E2").Value = 0.5771960150960151
$ws.Range("F2").Value = 0.8333333333333334
$ws.Range("G2").Value = 0.9698712042489004
$ws.Range("H2").Value = 0.5664023809523809
$ws.Range("I2").Value = 0.75
$ws.Range("J2").Value = 0.9872340425531916
$ws.Range("K2").Value = 0.612
$ws.Range("L2").Value = 0.9375
$ws.Range("M2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range("N2").Value = '[1 1 1 1 1 1 0 1 1 1 0 1 1 1 1 0 1 1 0 1 1 1 1 1]'

$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fa0044938b0>),
                (''model'',
                 XGBClassifier(base_score=None, booster=None, callbacks=None,
                               colsample_bylevel=None, colsample_bynode=None,
                               colsample_bytree=0.5, early_stopping_rounds=None,
                               enable_categorical=False, eval_metric=None,
                               feature_types=None, gamma=0.2, gpu_id=None,
                               grow_policy=None, importance_type=None,
                               interaction_constraints=None, learning_rate=0.01,
                               max_bin=None, max_cat_threshold=None,
                               max_cat_to_onehot=None, max_delta_step=None,
                               max_depth=7, max_leaves=None,
                               min_child_weight=None, missing=nan,
                               monotone_constraints=None, n_estimators=50,
                               n_jobs=None, num_parallel_tree=None,
                               predictor=None, random_state=42, ...))])'
$ws.Range("B3").Value = 0.6785880785880786
$ws.Range("C3").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fa0044dc580>, ''scaler'': RobustScaler(), ''model__subsample'': 0.5, ''model__n_estimators'': 50, ''model__max_depth'': 7, ''model__learning_rate'': 0.01, ''model__gamma'': 0.2, ''model__colsample_bytree'': 0.5}'
$ws.Range("D3").Value = 0.9671953648059105
$ws.Range("E3").Value = 0.5460047841047841
$ws.Range("F3").Value = 0.8421052631578948
$ws.Range("G3").Value = 0.9534032519674004
$ws.Range("H3").Value = 0.5796285714285714
$ws.Range("I3").Value = 0.7272727272727273
$ws.Range("J3").Value = 0.9828936170212766
$ws.Range("K3").Value = 0.5353333333333333
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]'
$ws.Range("N3").Value = '[1 1 1 1 1 0 1 1 1 1 1 0 1 1 1 1 1 1 1 1 1 1 1 1]'

$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fa004493490>),
                (''model'',
                 XGBClassifier(base_score=None, booster=None, callbacks=None,
                               colsample_bylevel=None, colsample_bynode=None,
                               colsample_bytree=0.8, early_stopping_rounds=None,
                               enable_categorical=False, eval_metric=None,
                               feature_types=None, gamma=0.2, gpu_id=None,
                               grow_policy=None, importance_type=None,
                               interaction_constraints=None, learning_rate=0.01,
                               max_bin=None, max_cat_threshold=None,
                               max_cat_to_onehot=None, max_delta_step=None,
                               max_depth=5, max_leaves=None,
                               min_child_weight=None, missing=nan,
                               monotone_constraints=None, n_estimators=100,
                               n_jobs=None, num_parallel_tree=None,
                               predictor=None, random_state=42, ...))])'
$ws.Range("B4").Value = 0.6686580086580087
$ws.Range("C4").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fa0044dc340>, ''scaler'': MinMaxScaler(), ''model__subsample'': 0.8, ''model__n_estimators'': 100, ''model__max_depth'': 5, ''model__learning_rate'': 0.01, ''model__gamma'': 0.2, ''model__colsample_bytree'': 0.8}'
$ws.Range("D4").Value = 0.9793729524970384
$ws.Range("E4").Value = 0.5626001887001887
$ws.Range("F4").Value = 0.7272727272727273
$ws.Range("G4").Value = 0.9731783971778825
$ws.Range("H4").Value = 0.5790111111111111
$ws.Range("I4").Value = 0.8571428571428571
$ws.Range("J4").Value = 0.9861333333333333
$ws.Range("K4").Value = 0.5723999999999999
$ws.Range("L4").Value = 0.631578947368421
$ws.Range("M4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range("N4").Value = '[0 1 1 0 0 1 0 1 1 0 0 0 0 1 1 1 0 1 1 0 1 1 1 1]'

$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fa0044dcca0>),
                (''model'',
                 XGBClassifier(base_score=None, booster=None, callbacks=None,
                               colsample_bylevel=None, colsample_bynode=None,
                               colsample_bytree=0.5, early_stopping_rounds=None,
                               enable_categorical=False, eval_metric=None,
                               feature_types=None, gamma=0, gpu_id=None,
                               grow_policy=None, importance_type=None,
                               interaction_constraints=None, learning_rate=0.01,
                               max_bin=None, max_cat_threshold=None,
                               max_cat_to_onehot=None, max_delta_step=None,
                               max_depth=3, max_leaves=None,
                               min_child_weight=None, missing=nan,
                               monotone_constraints=None, n_estimators=50,
                               n_jobs=None, num_parallel_tree=None,
                               predictor=None, random_state=42, ...))])'
$ws.Range("B5").Value = 0.7295787545787545
$ws.Range("C5").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fa004472070>, ''scaler'': RobustScaler(), ''model__subsample'': 0.5, ''model__n_estimators'': 50, ''model__max_depth'': 3, ''model__learning_rate'': 0.01, ''model__gamma'': 0, ''model__colsample_bytree'': 0.5}'
$ws.Range("D5").Value = 0.9724292358290221
$ws.Range("E5").Value = 0.6013865356865357
$ws.Range("F5").Value = 0.8125000000000001
$ws.Range("G5").Value = 0.9559223380217277
$ws.Range("H5").Value = 0.5739285714285715
$ws.Range("I5").Value = 0.7222222222222222
$ws.Range("J5").Value = 0.9918367346938777
$ws.Range("K5").Value = 0.6493333333333333
$ws.Range("L5").Value = 0.9285714285714286
$ws.Range("M5").Value = '[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]'
$ws.Range("N5").Value = '[0 1 1 1 1 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 0 1 1 1]'

$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fa004456910>),
                (''model'',
                 XGBClassifier(base_score=None, booster=None, callbacks=None,
                               colsample_bylevel=None, colsample_bynode=None,
                               colsample_bytree=0.5, early_stopping_rounds=None,
                               enable_categorical=False, eval_metric=None,
                               feature_types=None, gamma=0, gpu_id=None,
                               grow_policy=None, importance_type=None,
                               interaction_constraints=None, learning_rate=0.01,
                               max_bin=None, max_cat_threshold=None,
                               max_cat_to_onehot=None, max_delta_step=None,
                               max_depth=5, max_leaves=None,
                               min_child_weight=None, missing=nan,
                               monotone_constraints=None, n_estimators=50,
                               n_jobs=None, num_parallel_tree=None,
                               predictor=None, random_state=42, ...))])'
$ws.Range("B6").Value = 0.7458591408591408
$ws.Range("C6").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fa004420100>, ''scaler'': MinMaxScaler(), ''model__subsample'': 0.5, ''model__n_estimators'': 50, ''model__max_depth'': 5, ''model__learning_rate'': 0.01, ''model__gamma'': 0, ''model__colsample_bytree'': 0.5}'
$ws.Range("D6").Value = 0.968782605841421
$ws.Range("E6").Value = 0.6421301920301921
$ws.Range("F6").Value = 0.6451612903225806
$ws.Range("G6").Value = 0.9498296634007943
$ws.Range("H6").Value = 0.6088182539682541
$ws.Range("I6").Value = 0.5
$ws.Range("J6").Value = 0.9904999999999999
$ws.Range("K6").Value = 0.698
$ws.Range("L6").Value = 0.9090909090909091
$ws.Range("M6").Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]'
$ws.Range("N6").Value = '[1 1 1 1 1 1 1 0 1 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1]'
